# Refresh the cryptos list (Price / Volume(1h) columns) with the latest
# scraped figures. For cells whose new text looks like a plain number
# (e.g. "1.00", "8.00"), force text formatting first so Excel doesn't
# silently coerce the string into a numeric value (which would drop
# trailing zeros / reformat it) - then restore the "Normal" style so no
# stray number-format style is left attached to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.945.93'
$ws.Range("E2").Value = '  -3.15%  '

$ws.Range("D3").Value = '2.296.83'
$ws.Range("E3").Value = '  -3.53%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '534.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.05%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.59'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.20%  '

$ws.Range("E7").Value = '  +0.12%  '

$ws.Range("E8").Value = '  -1.16%  '

$ws.Range("D9").Value = '2.296.72'
$ws.Range("E9").Value = '  -3.50%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0996'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.42'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.77%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.148'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.51%  '

$ws.Range("E13").Value = '  -3.92%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.44'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.07%  '

$ws.Range("D15").Value = '2.708.22'
$ws.Range("E15").Value = '  -3.44%  '

$ws.Range("D16").Value = '57.924.80'
$ws.Range("E16").Value = '  -3.10%  '

$ws.Range("E17").Value = '  -4.19%  '

$ws.Range("D18").Value = '2.300.01'
$ws.Range("E18").Value = '  -4.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.50'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.57%  '

$ws.Range("E20").Value = '  -5.96%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '313.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.33%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.36'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.82%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.19%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.50'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.52%  '

$ws.Range("E25").Value = '  -4.28%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.994'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.54%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.76%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.28'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.47%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.78'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.46%  '

$ws.Range("E30").Value = '  -5.43%  '

$ws.Range("D31").Value = '0.0₃0716'
$ws.Range("E31").Value = '  -5.72%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.77'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.79%  '

$ws.Range("E33").Value = '  -4.62%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.380'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.42%  '

$ws.Range("E35").Value = '  -0.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.73'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.46%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.31%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.23'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.30%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.89'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.90%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '38.11'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.35%  '

$ws.Range("E41").Value = '  -6.00%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '140.33'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.48%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '287.55'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -9.94%  '

$ws.Range("E44").Value = '  -3.34%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0949'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.18%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0496'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.80%  '

$ws.Range("E47").Value = '  -2.95%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.04'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -8.27%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0210'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.56%  '

$ws.Range("E50").Value = '  -1.27%  '

$ws.Range("D51").Value = '0.0₆0205'
$ws.Range("E51").Value = '  +89.99%  '
